# Update Bmp2-Acvr2a NATMI ligand-receptor pair metrics
# (Ligand-expressing cells / Receptor-expressing cells count changed from 1 to 3 replicate,
# so all dependent expression/specificity/edge-weight metrics are recomputed)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 5.423951666666667
$ws.Cells.Item(2, 8).Value = 16.271855
$ws.Cells.Item(2, 9).Value = 0.4774188439413272
$ws.Cells.Item(2, 10).Value = 0.4774188439413271
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 17.13024733333333
$ws.Cells.Item(2, 14).Value = 51.390742
$ws.Cells.Item(2, 15).Value = 0.2959211466465044
$ws.Cells.Item(2, 16).Value = 0.2959211466465043
$ws.Cells.Item(2, 17).Value = 92.91363357404555
$ws.Cells.Item(2, 18).Value = 836.2227021664098
$ws.Cells.Item(2, 19).Value = 0.1412783317297661
$ws.Cells.Item(2, 20).Value = 0.141278331729766

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 5.423951666666667
$ws.Cells.Item(3, 8).Value = 16.271855
$ws.Cells.Item(3, 9).Value = 0.4774188439413272
$ws.Cells.Item(3, 10).Value = 0.4774188439413271
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 27.61090666666666
$ws.Cells.Item(3, 14).Value = 82.83272
$ws.Cells.Item(3, 15).Value = 0.4769721651858779
$ws.Cells.Item(3, 16).Value = 0.4769721651858778
$ws.Cells.Item(3, 17).Value = 149.7602232328444
$ws.Cells.Item(3, 18).Value = 1347.8420090956
$ws.Cells.Item(3, 19).Value = 0.2277154996952336
$ws.Cells.Item(3, 20).Value = 0.2277154996952335

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 5.423951666666667
$ws.Cells.Item(4, 8).Value = 16.271855
$ws.Cells.Item(4, 9).Value = 0.4774188439413272
$ws.Cells.Item(4, 10).Value = 0.4774188439413271
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 4.423514
$ws.Cells.Item(4, 14).Value = 13.270542
$ws.Cells.Item(4, 15).Value = 0.07641520344774541
$ws.Cells.Item(4, 16).Value = 0.0764152034477454
$ws.Cells.Item(4, 17).Value = 23.99292613282333
$ws.Cells.Item(4, 18).Value = 215.93633519541
$ws.Cells.Item(4, 19).Value = 0.03648205808956394
$ws.Cells.Item(4, 20).Value = 0.03648205808956392

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 5.423951666666667
$ws.Cells.Item(5, 8).Value = 16.271855
$ws.Cells.Item(5, 9).Value = 0.4774188439413272
$ws.Cells.Item(5, 10).Value = 0.4774188439413271
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 8.723210333333332
$ws.Cells.Item(5, 14).Value = 26.169631
$ws.Cells.Item(5, 15).Value = 0.1506914847198724
$ws.Cells.Item(5, 16).Value = 0.1506914847198724
$ws.Cells.Item(5, 17).Value = 47.31427122616721
$ws.Cells.Item(5, 18).Value = 425.828441035505
$ws.Cells.Item(5, 19).Value = 0.07194295442676367
$ws.Cells.Item(5, 20).Value = 0.07194295442676364

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 1.583504333333333
$ws.Cells.Item(6, 8).Value = 4.750513
$ws.Cells.Item(6, 9).Value = 0.1393808158066948
$ws.Cells.Item(6, 10).Value = 0.1393808158066948
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 17.13024733333333
$ws.Cells.Item(6, 14).Value = 51.390742
$ws.Cells.Item(6, 15).Value = 0.2959211466465044
$ws.Cells.Item(6, 16).Value = 0.2959211466465043
$ws.Cells.Item(6, 17).Value = 27.12582088340511
$ws.Cells.Item(6, 18).Value = 244.132387950646
$ws.Cells.Item(6, 19).Value = 0.04124573083404235
$ws.Cells.Item(6, 20).Value = 0.04124573083404234

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 1.583504333333333
$ws.Cells.Item(7, 8).Value = 4.750513
$ws.Cells.Item(7, 9).Value = 0.1393808158066948
$ws.Cells.Item(7, 10).Value = 0.1393808158066948
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 27.61090666666666
$ws.Cells.Item(7, 14).Value = 82.83272
$ws.Cells.Item(7, 15).Value = 0.4769721651858779
$ws.Cells.Item(7, 16).Value = 0.4769721651858778
$ws.Cells.Item(7, 17).Value = 43.72199035392888
$ws.Cells.Item(7, 18).Value = 393.4979131853599
$ws.Cells.Item(7, 19).Value = 0.06648076950069326
$ws.Cells.Item(7, 20).Value = 0.06648076950069323

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 1.583504333333333
$ws.Cells.Item(8, 8).Value = 4.750513
$ws.Cells.Item(8, 9).Value = 0.1393808158066948
$ws.Cells.Item(8, 10).Value = 0.1393808158066948
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 4.423514
$ws.Cells.Item(8, 14).Value = 13.270542
$ws.Cells.Item(8, 15).Value = 0.07641520344774541
$ws.Cells.Item(8, 16).Value = 0.0764152034477454
$ws.Cells.Item(8, 17).Value = 7.004653587560666
$ws.Cells.Item(8, 18).Value = 63.04188228804599
$ws.Cells.Item(8, 19).Value = 0.01065081339658131
$ws.Cells.Item(8, 20).Value = 0.01065081339658131

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 1.583504333333333
$ws.Cells.Item(9, 8).Value = 4.750513
$ws.Cells.Item(9, 9).Value = 0.1393808158066948
$ws.Cells.Item(9, 10).Value = 0.1393808158066948
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 8.723210333333332
$ws.Cells.Item(9, 14).Value = 26.169631
$ws.Cells.Item(9, 15).Value = 0.1506914847198724
$ws.Cells.Item(9, 16).Value = 0.1506914847198724
$ws.Cells.Item(9, 17).Value = 13.81324136341144
$ws.Cells.Item(9, 18).Value = 124.319172270703
$ws.Cells.Item(9, 19).Value = 0.02100350207537791
$ws.Cells.Item(9, 20).Value = 0.0210035020753779

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 0.6660723333333333
$ws.Cells.Item(10, 8).Value = 1.998217
$ws.Cells.Item(10, 9).Value = 0.05862800830537802
$ws.Cells.Item(10, 10).Value = 0.05862800830537802
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 17.13024733333333
$ws.Cells.Item(10, 14).Value = 51.390742
$ws.Cells.Item(10, 15).Value = 0.2959211466465044
$ws.Cells.Item(10, 16).Value = 0.2959211466465043
$ws.Cells.Item(10, 17).Value = 11.40998381189044
$ws.Cells.Item(10, 18).Value = 102.689854307014
$ws.Cells.Item(10, 19).Value = 0.01734926744332825
$ws.Cells.Item(10, 20).Value = 0.01734926744332824

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 0.6660723333333333
$ws.Cells.Item(11, 8).Value = 1.998217
$ws.Cells.Item(11, 9).Value = 0.05862800830537802
$ws.Cells.Item(11, 10).Value = 0.05862800830537802
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 27.61090666666666
$ws.Cells.Item(11, 14).Value = 82.83272
$ws.Cells.Item(11, 15).Value = 0.4769721651858779
$ws.Cells.Item(11, 16).Value = 0.4769721651858778
$ws.Cells.Item(11, 17).Value = 18.39086102891555
$ws.Cells.Item(11, 18).Value = 165.51774926024
$ws.Cells.Item(11, 19).Value = 0.02796392806195179
$ws.Cells.Item(11, 20).Value = 0.02796392806195178

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 0.6660723333333333
$ws.Cells.Item(12, 8).Value = 1.998217
$ws.Cells.Item(12, 9).Value = 0.05862800830537802
$ws.Cells.Item(12, 10).Value = 0.05862800830537802
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 4.423514
$ws.Cells.Item(12, 14).Value = 13.270542
$ws.Cells.Item(12, 15).Value = 0.07641520344774541
$ws.Cells.Item(12, 16).Value = 0.0764152034477454
$ws.Cells.Item(12, 17).Value = 2.946380291512666
$ws.Cells.Item(12, 18).Value = 26.517422623614
$ws.Cells.Item(12, 19).Value = 0.004480071182391569
$ws.Cells.Item(12, 20).Value = 0.004480071182391568

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 0.6660723333333333
$ws.Cells.Item(13, 8).Value = 1.998217
$ws.Cells.Item(13, 9).Value = 0.05862800830537802
$ws.Cells.Item(13, 10).Value = 0.05862800830537802
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 8.723210333333332
$ws.Cells.Item(13, 14).Value = 26.169631
$ws.Cells.Item(13, 15).Value = 0.1506914847198724
$ws.Cells.Item(13, 16).Value = 0.1506914847198724
$ws.Cells.Item(13, 17).Value = 5.810289060880777
$ws.Cells.Item(13, 18).Value = 52.29260154792699
$ws.Cells.Item(13, 19).Value = 0.008834741617706426
$ws.Cells.Item(13, 20).Value = 0.008834741617706424

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 3.687463666666666
$ws.Cells.Item(14, 8).Value = 11.062391
$ws.Cells.Item(14, 9).Value = 0.3245723319466
$ws.Cells.Item(14, 10).Value = 0.3245723319466
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 17.13024733333333
$ws.Cells.Item(14, 14).Value = 51.390742
$ws.Cells.Item(14, 15).Value = 0.2959211466465044
$ws.Cells.Item(14, 16).Value = 0.2959211466465043
$ws.Cells.Item(14, 17).Value = 63.16716464268022
$ws.Cells.Item(14, 18).Value = 568.5044817841219
$ws.Cells.Item(14, 19).Value = 0.0960478166393677
$ws.Cells.Item(14, 20).Value = 0.09604781663936768

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 3.687463666666666
$ws.Cells.Item(15, 8).Value = 11.062391
$ws.Cells.Item(15, 9).Value = 0.3245723319466
$ws.Cells.Item(15, 10).Value = 0.3245723319466
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 27.61090666666666
$ws.Cells.Item(15, 14).Value = 82.83272
$ws.Cells.Item(15, 15).Value = 0.4769721651858779
$ws.Cells.Item(15, 16).Value = 0.4769721651858778
$ws.Cells.Item(15, 17).Value = 101.8142151370578
$ws.Cells.Item(15, 18).Value = 916.32793623352
$ws.Cells.Item(15, 19).Value = 0.1548119679279993
$ws.Cells.Item(15, 20).Value = 0.1548119679279993

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 3.687463666666666
$ws.Cells.Item(16, 8).Value = 11.062391
$ws.Cells.Item(16, 9).Value = 0.3245723319466
$ws.Cells.Item(16, 10).Value = 0.3245723319466
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 4.423514
$ws.Cells.Item(16, 14).Value = 13.270542
$ws.Cells.Item(16, 15).Value = 0.07641520344774541
$ws.Cells.Item(16, 16).Value = 0.0764152034477454
$ws.Cells.Item(16, 17).Value = 16.31154715399133
$ws.Cells.Item(16, 18).Value = 146.803924385922
$ws.Cells.Item(16, 19).Value = 0.02480226077920859
$ws.Cells.Item(16, 20).Value = 0.02480226077920859

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 3.687463666666666
$ws.Cells.Item(17, 8).Value = 11.062391
$ws.Cells.Item(17, 9).Value = 0.3245723319466
$ws.Cells.Item(17, 10).Value = 0.3245723319466
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 8.723210333333332
$ws.Cells.Item(17, 14).Value = 26.169631
$ws.Cells.Item(17, 15).Value = 0.1506914847198724
$ws.Cells.Item(17, 16).Value = 0.1506914847198724
$ws.Cells.Item(17, 17).Value = 32.16652116085788
$ws.Cells.Item(17, 18).Value = 289.498690447721
$ws.Cells.Item(17, 19).Value = 0.04891028660002443
$ws.Cells.Item(17, 20).Value = 0.04891028660002442

